$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.792.57"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "2.949.84"
$ws.Range("E3").Value = "  +3.74%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "352.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.65%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "112.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.70%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.564"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.07%  "
$ws.Range("E9").Value = "  +1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.11%  "
$ws.Range("E11").Value = "  +4.90%  "
$ws.Range("E12").Value = "  +1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "20.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.87"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "3.410.48"
$ws.Range("E15").Value = "  +4.04%  "
$ws.Range("D16").Value = "2.934.95"
$ws.Range("E16").Value = "  +2.90%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.987"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.01%  "
$ws.Range("D18").Value = "51.890.59"
$ws.Range("E18").Value = "  -0.05%  "
$ws.Range("E19").Value = "  +0.03%  "
$ws.Range("E20").Value = "  -2.70%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.47"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.16%  "
$ws.Range("D22").Value = "0.0₃0985"
$ws.Range("E22").Value = "  +1.16%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "71.33"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "269.42"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.80"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.39%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.180"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +9.52%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "27.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.08%  "
$ws.Range("E28").Value = "  +0.01%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.42"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +17.44%  "
$ws.Range("E30").Value = "  +25.52%  "
$ws.Range("E32").Value = "  +0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "37.49"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "52.88"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.21%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "18.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.84%  "
$ws.Range("E40").Value = "  +1.50%  "
$ws.Range("E41").Value = "  +5.48%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.117"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.43%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "23.54"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +6.05%  "
$ws.Range("E44").Value = "  -1.12%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.53"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.34%  "
$ws.Range("E46").Value = "  -0.30%  "
$ws.Range("D47").Value = "2.168.68"
$ws.Range("E47").Value = "  -0.38%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "112.98"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -8.27%  "
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0341"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +9.29%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.935"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.84%  "
